$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rotate labels C1,D1,E1 -> new C1="prediction", D1="rejection-f", E1="max"
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data rows 2-9: column C becomes the species prediction text (same as column D),
# column D stays the same, column E becomes numeric 1
for ($r = 2; $r -le 9; $r++) {
    $species = $ws.Cells.Item($r, 4).Text
    $ws.Cells.Item($r, 3).Value = $species
    $ws.Cells.Item($r, 5).Value = 1
}
